$wb = $excel.ActiveWorkbook

# Excel border/line-style/weight constants
$edgeTop    = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop
$edgeBottom = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom
$edgeRight  = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight
$continuous = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$thin       = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

function Set-TopBottomBorder($rng) {
    $rng.ClearFormats()
    $rng.Borders.Item($edgeTop).LineStyle = $continuous
    $rng.Borders.Item($edgeTop).Weight = $thin
    $rng.Borders.Item($edgeBottom).LineStyle = $continuous
    $rng.Borders.Item($edgeBottom).Weight = $thin
}

function Set-TopBottomRightBorder($rng) {
    $rng.ClearFormats()
    $rng.Borders.Item($edgeTop).LineStyle = $continuous
    $rng.Borders.Item($edgeTop).Weight = $thin
    $rng.Borders.Item($edgeBottom).LineStyle = $continuous
    $rng.Borders.Item($edgeBottom).Weight = $thin
    $rng.Borders.Item($edgeRight).LineStyle = $continuous
    $rng.Borders.Item($edgeRight).Weight = $thin
}

# ---------- Sheet 1: quality_comparison ----------
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder      $ws1.Range("C1")
Set-TopBottomRightBorder $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# ---------- Sheet 2: computational_comparison ----------
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder      $ws2.Range("C1")
Set-TopBottomRightBorder $ws2.Range("D1")
Set-TopBottomBorder      $ws2.Range("F1")
Set-TopBottomRightBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
